# Applies updated market-price / profit figures to the Leve profit tables
# across all 8 crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 862.5
$ws.Range("I18").Value = 916.6667
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 916.6667
$ws.Range("L18").Value = 700
$ws.Range("M18").Value = -632.6667
$ws.Range("N18").Value = -1268

$ws.Range("H70").Value = 1046.6428
$ws.Range("I70").Value = 1104.5454
$ws.Range("J70").Value = 834.3333
$ws.Range("K70").Value = 3313.6362
$ws.Range("L70").Value = 2502.9999
$ws.Range("M70").Value = -3043.6362
$ws.Range("N70").Value = -3042.9999

$ws.Range("H73").Value = 1046.6428
$ws.Range("I73").Value = 1104.5454
$ws.Range("J73").Value = 834.3333
$ws.Range("K73").Value = 3313.6362
$ws.Range("L73").Value = 2502.9999
$ws.Range("M73").Value = -2377.6362
$ws.Range("N73").Value = -4374.9999

$ws.Range("H80").Value = 2657599
$ws.Range("I80").Value = 1287.875
$ws.Range("J80").Value = 3669527
$ws.Range("K80").Value = 3863.625
$ws.Range("L80").Value = 11008581
$ws.Range("M80").Value = -2865.625
$ws.Range("N80").Value = -11010577

$ws.Range("H83").Value = 2657599
$ws.Range("I83").Value = 1287.875
$ws.Range("J83").Value = 3669527
$ws.Range("K83").Value = 11590.875
$ws.Range("L83").Value = 33025743
$ws.Range("M83").Value = -6598.875
$ws.Range("N83").Value = -33035727

$ws.Range("H86").Value = 7698.875
$ws.Range("I86").Value = 1520.2
$ws.Range("K86").Value = 1520.2
$ws.Range("M86").Value = -397.2

$ws.Range("H89").Value = 7698.875
$ws.Range("I89").Value = 1520.2
$ws.Range("K89").Value = 7601
$ws.Range("M89").Value = -1985

$ws.Range("H127").Value = 1208.0714
$ws.Range("I127").Value = 758.125
$ws.Range("K127").Value = 2274.375
$ws.Range("M127").Value = 2685.625

$ws.Range("H129").Value = 795.96906
$ws.Range("J129").Value = 827.93256
$ws.Range("L129").Value = 2483.79768
$ws.Range("N129").Value = -12483.79768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 30376.5
$ws.Range("J43").Value = 30376.5
$ws.Range("L43").Value = 30376.5
$ws.Range("N43").Value = -31002.5

$ws.Range("H61").Value = 3554.0952
$ws.Range("I61").Value = 2188.7778
$ws.Range("J61").Value = 4578.0835
$ws.Range("K61").Value = 2188.7778
$ws.Range("L61").Value = 4578.0835
$ws.Range("M61").Value = -1976.7778
$ws.Range("N61").Value = -5002.0835

$ws.Range("H136").Value = 3554.0952
$ws.Range("I136").Value = 2188.7778
$ws.Range("J136").Value = 4578.0835
$ws.Range("K136").Value = 6566.3334
$ws.Range("L136").Value = 13734.2505
$ws.Range("M136").Value = -4016.3334
$ws.Range("N136").Value = -18834.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2033.6875
$ws.Range("I20").Value = 2236.5833
$ws.Range("J20").Value = 1425
$ws.Range("K20").Value = 2236.5833
$ws.Range("L20").Value = 1425
$ws.Range("M20").Value = -1989.5833
$ws.Range("N20").Value = -1919

$ws.Range("H41").Value = 45992.5
$ws.Range("J41").Value = 45992.5
$ws.Range("L41").Value = 45992.5
$ws.Range("N41").Value = -46768.5

$ws.Range("H48").Value = 119995
$ws.Range("J48").Value = 119995
$ws.Range("L48").Value = 119995
$ws.Range("N48").Value = -120825

$ws.Range("H119").Value = 20761
$ws.Range("J119").Value = 20761
$ws.Range("L119").Value = 20761
$ws.Range("N119").Value = -30437

$ws.Range("H120").Value = 37630
$ws.Range("J120").Value = 37630
$ws.Range("L120").Value = 37630
$ws.Range("N120").Value = -47306

$ws.Range("H134").Value = 3219.75
$ws.Range("I134").Value = 3406.12
$ws.Range("J134").Value = 1666.6666
$ws.Range("K134").Value = 10218.36
$ws.Range("L134").Value = 4999.9998
$ws.Range("M134").Value = -7683.360000000001
$ws.Range("N134").Value = -10069.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 102161660
$ws.Range("I6").Value = 25742500
$ws.Range("K6").Value = 25742500
$ws.Range("M6").Value = -25742387

$ws.Range("H18").Value = 35634.5
$ws.Range("J18").Value = 35634.5
$ws.Range("L18").Value = 35634.5
$ws.Range("N18").Value = -36094.5

$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150

$ws.Range("H122").Value = 2671.8572
$ws.Range("I122").Value = 2671.8572
$ws.Range("K122").Value = 8015.571599999999
$ws.Range("M122").Value = -5565.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 50150
$ws.Range("J107").Value = 300
$ws.Range("L107").Value = 900
$ws.Range("N107").Value = -4740

$ws.Range("H109").Value = 3885.1667
$ws.Range("I109").Value = 775.5714
$ws.Range("J109").Value = 5864
$ws.Range("K109").Value = 2326.7142
$ws.Range("L109").Value = 17592
$ws.Range("M109").Value = -1286.7142
$ws.Range("N109").Value = -19672

$ws.Range("H114").Value = 176.27272
$ws.Range("I114").Value = 128.28572
$ws.Range("J114").Value = 260.25
$ws.Range("K114").Value = 384.85716
$ws.Range("L114").Value = 780.75
$ws.Range("M114").Value = 2869.14284
$ws.Range("N114").Value = -7288.75

$ws.Range("H131").Value = 774.39
$ws.Range("I131").Value = 445.8
$ws.Range("J131").Value = 791.6842
$ws.Range("K131").Value = 1337.4
$ws.Range("L131").Value = 2375.0526
$ws.Range("M131").Value = 3702.6
$ws.Range("N131").Value = -12455.0526

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 103500
$ws.Range("J20").Value = 103500
$ws.Range("L20").Value = 103500
$ws.Range("N20").Value = -103990

$ws.Range("H70").Value = 4475028.5
$ws.Range("I70").Value = 22119.6
$ws.Range("K70").Value = 22119.6
$ws.Range("M70").Value = -21849.6

$ws.Range("H73").Value = 4475028.5
$ws.Range("I73").Value = 22119.6
$ws.Range("K73").Value = 22119.6
$ws.Range("M73").Value = -21183.6

$ws.Range("H122").Value = 2085.85
$ws.Range("I122").Value = 1758.5
$ws.Range("J122").Value = 2576.875
$ws.Range("K122").Value = 5275.5
$ws.Range("L122").Value = 7730.625
$ws.Range("M122").Value = -2825.5
$ws.Range("N122").Value = -12630.625

$ws.Range("H132").Value = 29939.422
$ws.Range("I132").Value = 3675.0715
$ws.Range("J132").Value = 103479.6
$ws.Range("K132").Value = 11025.2145
$ws.Range("L132").Value = 310438.8
$ws.Range("M132").Value = -8495.2145
$ws.Range("N132").Value = -315498.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2900

$ws.Range("H132").Value = 2040.4375
$ws.Range("I132").Value = 1249.909
$ws.Range("K132").Value = 3749.727
$ws.Range("M132").Value = -1219.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 40301
$ws.Range("J82").Value = 40301
$ws.Range("L82").Value = 40301
$ws.Range("N82").Value = -41067

$ws.Range("H85").Value = 40301
$ws.Range("J85").Value = 40301
$ws.Range("L85").Value = 40301
$ws.Range("N85").Value = -42953

$ws.Range("H132").Value = 1031.238
$ws.Range("I132").Value = 613.5789
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 1840.7367
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = 689.2633000000001
$ws.Range("N132").Value = -20057

$ws.Range("H136").Value = 16950278
$ws.Range("I136").Value = 22223144
$ws.Range("J136").Value = 1783.6428
$ws.Range("K136").Value = 66669432
$ws.Range("L136").Value = 5350.928400000001
$ws.Range("M136").Value = -66666882
$ws.Range("N136").Value = -10450.9284
